# Updated cryptos list on Mon Sep 30 23:27:09 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.495.66"
$ws.Range("E2").Value = "  -3.22%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.605.70"
$ws.Range("E3").Value = "  -1.86%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.26"
$ws.Range("E5").Value = "  -4.22%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.05"
$ws.Range("E6").Value = "  -2.78%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.619"
$ws.Range("E8").Value = "  -3.40%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.604.18"
$ws.Range("E9").Value = "  -1.80%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.115"
$ws.Range("E10").Value = "  -8.09%  "

# Row 11 - Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.79"
$ws.Range("E11").Value = "  -0.82%  "

# Row 12 - Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.377"
$ws.Range("E12").Value = "  -5.26%  "

# Row 13 - TRON
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.157"
$ws.Range("E13").Value = "  -0.03%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.91"
$ws.Range("E14").Value = "  -3.92%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.074.38"
$ws.Range("E15").Value = "  -1.84%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -7.72%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "63.389.20"
$ws.Range("E17").Value = "  -3.24%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.605.37"
$ws.Range("E18").Value = "  -2.25%  "

# Row 19 - Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.89"
$ws.Range("E19").Value = "  -4.77%  "

# Row 20 - Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.46"
$ws.Range("E20").Value = "  +0.43%  "

# Row 21 - Polkadot
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.47"
$ws.Range("E21").Value = "  -6.24%  "

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "339.49"
$ws.Range("E22").Value = "  -3.83%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.09%  "

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.01"
$ws.Range("E24").Value = "  -3.70%  "

# Row 25 - SuiNetwork
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.79"
$ws.Range("E25").Value = "  +1.32%  "

# Row 26 - PEPE
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000105"
$ws.Range("E26").Value = "  -6.74%  "

# Row 27 - Bittensor->InternetComputer(DFINITY)
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("E27").Value = "  -6.25%  "

# Row 28 - InternetComputer(DFINITY)->Bittensor
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "576.56"
$ws.Range("E28").Value = "  +3.13%  "

# Row 29 - Fetch.AI
$ws.Range("E29").Value = "  -4.68%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  -0.08%  "

# Row 31 - Kaspa
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.160"
$ws.Range("E31").Value = "  -2.06%  "

# Row 32 - Aptos
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.76"
$ws.Range("E32").Value = "  -3.93%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  -4.66%  "

# Row 34 - ImmutableX
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.71"
$ws.Range("E34").Value = "  -5.23%  "

# Row 35 - RenderToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.51"
$ws.Range("E35").Value = "  -2.22%  "

# Row 36 - NEARProtocol
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.33"
$ws.Range("E36").Value = "  -2.45%  "

# Row 37 - PolygonEcosystemToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.399"
$ws.Range("E37").Value = "  -5.29%  "

# Row 38 - FirstDigitalUSD
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.04%  "

# Row 39 - EthereumClassic
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.57"
$ws.Range("E39").Value = "  -4.37%  "

# Row 40 - Monero
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.56"
$ws.Range("E40").Value = "  +0.63%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  -5.27%  "

# Row 42 - USDe
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.01%  "

# Row 43 - OKB
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.56"
$ws.Range("E43").Value = "  -3.00%  "

# Row 44 - dogwifhat
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.47"
$ws.Range("E44").Value = "  +0.52%  "

# Row 45 - Aave
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "156.99"
$ws.Range("E45").Value = "  -2.47%  "

# Row 46 - InjectiveProtocol
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.34"
$ws.Range("E46").Value = "  +0.53%  "

# Row 47 - Filecoin
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.83"
$ws.Range("E47").Value = "  -6.17%  "

# Row 48 - Hedera
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0580"
$ws.Range("E48").Value = "  -5.73%  "

# Row 49 - Mantle
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.626"
$ws.Range("E49").Value = "  -2.59%  "

# Row 50 - Stellar
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0993"
$ws.Range("E50").Value = "  -2.19%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  -5.17%  "
